$wb = $excel.ActiveWorkbook

# --- Remove header-row (A1:N1) formatting on every sheet ---
foreach ($ws in $wb.Worksheets) {
    $ws.Range("A1:N1").ClearFormats()
}

# --- Update computed price/profit columns (H:N) with refreshed market data ---

# Sheet 1 = ALC
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 8).Value = 317.41666  # H2
$ws.Cells.Item(2, 9).Value = 301.7  # I2
$ws.Cells.Item(2, 10).Value = 396  # J2
$ws.Cells.Item(2, 11).Value = 301.7  # K2
$ws.Cells.Item(2, 12).Value = 396  # L2
$ws.Cells.Item(2, 13).Value = -188.7  # M2
$ws.Cells.Item(2, 14).Value = -622  # N2
$ws.Cells.Item(19, 8).Value = 755.9167  # H19
$ws.Cells.Item(19, 9).Value = 750  # I19
$ws.Cells.Item(19, 10).Value = 761.8333  # J19
$ws.Cells.Item(19, 11).Value = 750  # K19
$ws.Cells.Item(19, 12).Value = 761.8333  # L19
$ws.Cells.Item(19, 13).Value = -575  # M19
$ws.Cells.Item(19, 14).Value = -1111.8333  # N19
$ws.Cells.Item(64, 8).Value = 5689.6665  # H64
$ws.Cells.Item(64, 9).Value = 5669  # I64
$ws.Cells.Item(64, 10).Value = 5700  # J64
$ws.Cells.Item(64, 11).Value = 5669  # K64
$ws.Cells.Item(64, 12).Value = 5700  # L64
$ws.Cells.Item(64, 13).Value = -5421  # M64
$ws.Cells.Item(64, 14).Value = -6196  # N64
$ws.Cells.Item(67, 8).Value = 5689.6665  # H67
$ws.Cells.Item(67, 9).Value = 5669  # I67
$ws.Cells.Item(67, 10).Value = 5700  # J67
$ws.Cells.Item(67, 11).Value = 5669  # K67
$ws.Cells.Item(67, 12).Value = 5700  # L67
$ws.Cells.Item(67, 13).Value = -4811  # M67
$ws.Cells.Item(67, 14).Value = -7416  # N67
$ws.Cells.Item(70, 8).Value = 6716.9546  # H70
$ws.Cells.Item(70, 10).Value = 8426.333  # J70
$ws.Cells.Item(70, 12).Value = 25278.999  # L70
$ws.Cells.Item(70, 14).Value = -25818.999  # N70
$ws.Cells.Item(73, 8).Value = 6716.9546  # H73
$ws.Cells.Item(73, 10).Value = 8426.333  # J73
$ws.Cells.Item(73, 12).Value = 25278.999  # L73
$ws.Cells.Item(73, 14).Value = -27150.999  # N73
$ws.Cells.Item(76, 8).Value = 4436.5  # H76
$ws.Cells.Item(76, 9).Value = 4199.5  # I76
$ws.Cells.Item(76, 11).Value = 4199.5  # K76
$ws.Cells.Item(76, 13).Value = -3884.5  # M76
$ws.Cells.Item(79, 8).Value = 4436.5  # H79
$ws.Cells.Item(79, 9).Value = 4199.5  # I79
$ws.Cells.Item(79, 11).Value = 4199.5  # K79
$ws.Cells.Item(79, 13).Value = -3107.5  # M79
$ws.Cells.Item(132, 8).Value = 2455.6  # H132
$ws.Cells.Item(132, 9).Value = 2269.0527  # I132
$ws.Cells.Item(132, 11).Value = 6807.158100000001  # K132
$ws.Cells.Item(132, 13).Value = -4277.158100000001  # M132

# Sheet 2 = ARM
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(25, 8).Value = 2574.75  # H25
$ws.Cells.Item(25, 10).Value = 2599.6667  # J25
$ws.Cells.Item(25, 12).Value = 2599.6667  # L25
$ws.Cells.Item(25, 14).Value = -3403.6667  # N25
$ws.Cells.Item(34, 8).Value = 0  # H34
$ws.Cells.Item(34, 10).Value = 0  # J34
$ws.Cells.Item(34, 12).Value = 0  # L34
$ws.Cells.Item(34, 14).ClearContents()  # N34
$ws.Cells.Item(41, 8).Value = 3071.2  # H41
$ws.Cells.Item(41, 9).Value = 2964  # I41
$ws.Cells.Item(41, 11).Value = 2964  # K41
$ws.Cells.Item(41, 13).Value = -2550  # M41
$ws.Cells.Item(45, 8).Value = 2986.5  # H45
$ws.Cells.Item(45, 9).Value = 2803.5454  # I45
$ws.Cells.Item(45, 10).Value = 4999  # J45
$ws.Cells.Item(45, 11).Value = 2803.5454  # K45
$ws.Cells.Item(45, 12).Value = 4999  # L45
$ws.Cells.Item(45, 13).Value = -2426.5454  # M45
$ws.Cells.Item(45, 14).Value = -5753  # N45
$ws.Cells.Item(61, 8).Value = 2811.6775  # H61
$ws.Cells.Item(61, 9).Value = 1962.9048  # I61
$ws.Cells.Item(61, 10).Value = 4594.1  # J61
$ws.Cells.Item(61, 11).Value = 1962.9048  # K61
$ws.Cells.Item(61, 12).Value = 4594.1  # L61
$ws.Cells.Item(61, 13).Value = -1750.9048  # M61
$ws.Cells.Item(61, 14).Value = -5018.1  # N61
$ws.Cells.Item(94, 8).Value = 50000  # H94
$ws.Cells.Item(94, 10).Value = 50000  # J94
$ws.Cells.Item(94, 12).Value = 50000  # L94
$ws.Cells.Item(94, 14).Value = -51802  # N94
$ws.Cells.Item(122, 8).Value = 1664.409  # H122
$ws.Cells.Item(122, 9).Value = 1283.7858  # I122
$ws.Cells.Item(122, 11).Value = 3851.3574  # K122
$ws.Cells.Item(122, 13).Value = -1401.3574  # M122
$ws.Cells.Item(132, 8).Value = 2191.1333  # H132
$ws.Cells.Item(132, 9).Value = 1898.2593  # I132
$ws.Cells.Item(132, 11).Value = 5694.7779  # K132
$ws.Cells.Item(132, 13).Value = -3164.7779  # M132
$ws.Cells.Item(134, 8).Value = 100142.664  # H134
$ws.Cells.Item(134, 10).Value = 100142.664  # J134
$ws.Cells.Item(134, 12).Value = 100142.664  # L134
$ws.Cells.Item(134, 14).Value = -110282.664  # N134
$ws.Cells.Item(136, 8).Value = 2811.6775  # H136
$ws.Cells.Item(136, 9).Value = 1962.9048  # I136
$ws.Cells.Item(136, 10).Value = 4594.1  # J136
$ws.Cells.Item(136, 11).Value = 5888.7144  # K136
$ws.Cells.Item(136, 12).Value = 13782.3  # L136
$ws.Cells.Item(136, 13).Value = -3338.7144  # M136
$ws.Cells.Item(136, 14).Value = -18882.3  # N136

# Sheet 3 = BSM
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(75, 8).Value = 6348.4546  # H75
$ws.Cells.Item(75, 9).Value = 6683.3  # I75
$ws.Cells.Item(75, 11).Value = 6683.3  # K75
$ws.Cells.Item(75, 13).Value = -5747.3  # M75
$ws.Cells.Item(76, 8).Value = 32500  # H76
$ws.Cells.Item(76, 9).Value = 20000  # I76
$ws.Cells.Item(76, 11).Value = 20000  # K76
$ws.Cells.Item(76, 13).Value = -19685  # M76
$ws.Cells.Item(78, 8).Value = 6348.4546  # H78
$ws.Cells.Item(78, 9).Value = 6683.3  # I78
$ws.Cells.Item(78, 11).Value = 20049.9  # K78
$ws.Cells.Item(78, 13).Value = -15369.9  # M78
$ws.Cells.Item(79, 8).Value = 32500  # H79
$ws.Cells.Item(79, 9).Value = 20000  # I79
$ws.Cells.Item(79, 11).Value = 20000  # K79
$ws.Cells.Item(79, 13).Value = -18908  # M79
$ws.Cells.Item(96, 8).Value = 15665  # H96
$ws.Cells.Item(96, 9).Value = 15665  # I96
$ws.Cells.Item(96, 11).Value = 15665  # K96
$ws.Cells.Item(96, 13).Value = -12919  # M96

# Sheet 4 = CRP
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(7, 8).Value = 104.9  # H7
$ws.Cells.Item(7, 9).Value = 118.5  # I7
$ws.Cells.Item(7, 10).Value = 84.5  # J7
$ws.Cells.Item(7, 11).Value = 118.5  # K7
$ws.Cells.Item(7, 12).Value = 84.5  # L7
$ws.Cells.Item(7, 13).Value = -5.5  # M7
$ws.Cells.Item(7, 14).Value = -310.5  # N7
$ws.Cells.Item(31, 8).Value = 1548.8292  # H31
$ws.Cells.Item(31, 9).Value = 881.43475  # I31
$ws.Cells.Item(31, 10).Value = 2401.611  # J31
$ws.Cells.Item(31, 11).Value = 881.43475  # K31
$ws.Cells.Item(31, 12).Value = 2401.611  # L31
$ws.Cells.Item(31, 13).Value = -586.43475  # M31
$ws.Cells.Item(31, 14).Value = -2991.611  # N31
$ws.Cells.Item(34, 8).Value = 1548.8292  # H34
$ws.Cells.Item(34, 9).Value = 881.43475  # I34
$ws.Cells.Item(34, 10).Value = 2401.611  # J34
$ws.Cells.Item(34, 11).Value = 881.43475  # K34
$ws.Cells.Item(34, 12).Value = 2401.611  # L34
$ws.Cells.Item(34, 13).Value = -679.43475  # M34
$ws.Cells.Item(34, 14).Value = -2805.611  # N34
$ws.Cells.Item(58, 8).Value = 2863.125  # H58
$ws.Cells.Item(58, 9).Value = 3104.3333  # I58
$ws.Cells.Item(58, 10).Value = 2139.5  # J58
$ws.Cells.Item(58, 11).Value = 3104.3333  # K58
$ws.Cells.Item(58, 12).Value = 2139.5  # L58
$ws.Cells.Item(58, 13).Value = -2901.3333  # M58
$ws.Cells.Item(58, 14).Value = -2545.5  # N58
$ws.Cells.Item(93, 8).Value = 23215  # H93
$ws.Cells.Item(93, 9).Value = 23215  # I93
$ws.Cells.Item(93, 11).Value = 23215  # K93
$ws.Cells.Item(93, 13).Value = -21343  # M93
$ws.Cells.Item(99, 8).Value = 1588.8  # H99
$ws.Cells.Item(99, 9).Value = 1479.0834  # I99
$ws.Cells.Item(99, 11).Value = 1479.0834  # K99
$ws.Cells.Item(99, 13).Value = 18.91660000000002  # M99
$ws.Cells.Item(103, 8).Value = 23599.2  # H103
$ws.Cells.Item(103, 9).Value = 11999.25  # I103
$ws.Cells.Item(103, 11).Value = 11999.25  # K103
$ws.Cells.Item(103, 13).Value = -10827.25  # M103
$ws.Cells.Item(126, 8).Value = 1588.8  # H126
$ws.Cells.Item(126, 9).Value = 1479.0834  # I126
$ws.Cells.Item(126, 11).Value = 4437.2502  # K126
$ws.Cells.Item(126, 13).Value = -1967.2502  # M126
$ws.Cells.Item(136, 8).Value = 2863.125  # H136
$ws.Cells.Item(136, 9).Value = 3104.3333  # I136
$ws.Cells.Item(136, 10).Value = 2139.5  # J136
$ws.Cells.Item(136, 11).Value = 9312.999899999999  # K136
$ws.Cells.Item(136, 12).Value = 6418.5  # L136
$ws.Cells.Item(136, 13).Value = -6762.999899999999  # M136
$ws.Cells.Item(136, 14).Value = -11518.5  # N136

# Sheet 5 = CUL
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(10, 8).Value = 552.9  # H10
$ws.Cells.Item(10, 9).Value = 58.88889  # I10
$ws.Cells.Item(10, 11).Value = 176.66667  # K10
$ws.Cells.Item(10, 13).Value = -37.66667000000001  # M10
$ws.Cells.Item(20, 8).Value = 9450  # H20
$ws.Cells.Item(20, 10).Value = 9450  # J20
$ws.Cells.Item(20, 12).Value = 28350  # L20
$ws.Cells.Item(20, 14).Value = -28804  # N20
$ws.Cells.Item(21, 8).Value = 3249.6667  # H21
$ws.Cells.Item(21, 9).Value = 499.5  # I21
$ws.Cells.Item(21, 10).Value = 4624.75  # J21
$ws.Cells.Item(21, 11).Value = 1498.5  # K21
$ws.Cells.Item(21, 12).Value = 13874.25  # L21
$ws.Cells.Item(21, 13).Value = -1325.5  # M21
$ws.Cells.Item(21, 14).Value = -14220.25  # N21
$ws.Cells.Item(34, 8).Value = 5848.8335  # H34
$ws.Cells.Item(34, 10).Value = 5848.8335  # J34
$ws.Cells.Item(34, 12).Value = 17546.5005  # L34
$ws.Cells.Item(34, 14).Value = -17714.5005  # N34
$ws.Cells.Item(37, 8).Value = 150000  # H37
$ws.Cells.Item(37, 10).Value = 150000  # J37
$ws.Cells.Item(37, 12).Value = 450000  # L37
$ws.Cells.Item(37, 14).Value = -450224  # N37
$ws.Cells.Item(38, 8).Value = 130.09525  # H38
$ws.Cells.Item(38, 9).Value = 81.083336  # I38
$ws.Cells.Item(38, 11).Value = 243.250008  # K38
$ws.Cells.Item(38, 13).Value = 103.749992  # M38
$ws.Cells.Item(60, 8).Value = 933.6667  # H60
$ws.Cells.Item(60, 9).Value = 1172.8572  # I60
$ws.Cells.Item(60, 11).Value = 3518.5716  # K60
$ws.Cells.Item(60, 13).Value = -3267.5716  # M60
$ws.Cells.Item(70, 8).Value = 6094.2856  # H70
$ws.Cells.Item(70, 9).Value = 3932  # I70
$ws.Cells.Item(70, 11).Value = 11796  # K70
$ws.Cells.Item(70, 13).Value = -11481  # M70
$ws.Cells.Item(73, 8).Value = 6094.2856  # H73
$ws.Cells.Item(73, 9).Value = 3932  # I73
$ws.Cells.Item(73, 11).Value = 11796  # K73
$ws.Cells.Item(73, 13).Value = -10704  # M73
$ws.Cells.Item(107, 8).Value = 1192.6904  # H107
$ws.Cells.Item(107, 10).Value = 1348.9722  # J107
$ws.Cells.Item(107, 12).Value = 4046.9166  # L107
$ws.Cells.Item(107, 14).Value = -7886.9166  # N107
$ws.Cells.Item(116, 8).Value = 3450  # H116
$ws.Cells.Item(116, 9).Value = 3000  # I116
$ws.Cells.Item(116, 10).Value = 3900  # J116
$ws.Cells.Item(116, 11).Value = 9000  # K116
$ws.Cells.Item(116, 12).Value = 11700  # L116
$ws.Cells.Item(116, 13).Value = -5558  # M116
$ws.Cells.Item(116, 14).Value = -18584  # N116
$ws.Cells.Item(139, 8).Value = 18659.6  # H139
$ws.Cells.Item(139, 9).Value = 1465  # I139
$ws.Cells.Item(139, 11).Value = 4395  # K139
$ws.Cells.Item(139, 13).Value = 745  # M139

# Sheet 6 = GSM
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(70, 8).Value = 6714.636  # H70
$ws.Cells.Item(70, 9).Value = 6409.143  # I70
$ws.Cells.Item(70, 10).Value = 7249.25  # J70
$ws.Cells.Item(70, 11).Value = 6409.143  # K70
$ws.Cells.Item(70, 12).Value = 7249.25  # L70
$ws.Cells.Item(70, 13).Value = -6139.143  # M70
$ws.Cells.Item(70, 14).Value = -7789.25  # N70
$ws.Cells.Item(73, 8).Value = 6714.636  # H73
$ws.Cells.Item(73, 9).Value = 6409.143  # I73
$ws.Cells.Item(73, 10).Value = 7249.25  # J73
$ws.Cells.Item(73, 11).Value = 6409.143  # K73
$ws.Cells.Item(73, 12).Value = 7249.25  # L73
$ws.Cells.Item(73, 13).Value = -5473.143  # M73
$ws.Cells.Item(73, 14).Value = -9121.25  # N73
$ws.Cells.Item(80, 8).Value = 5333.3335  # H80
$ws.Cells.Item(80, 9).Value = 0  # I80
$ws.Cells.Item(80, 11).Value = 0  # K80
$ws.Cells.Item(80, 13).ClearContents()  # M80
$ws.Cells.Item(83, 8).Value = 5333.3335  # H83
$ws.Cells.Item(83, 9).Value = 0  # I83
$ws.Cells.Item(83, 11).Value = 0  # K83
$ws.Cells.Item(83, 13).ClearContents()  # M83
$ws.Cells.Item(87, 8).Value = 76499.75  # H87
$ws.Cells.Item(87, 10).Value = 76499.75  # J87
$ws.Cells.Item(87, 12).Value = 76499.75  # L87
$ws.Cells.Item(87, 14).Value = -78995.75  # N87
$ws.Cells.Item(90, 8).Value = 76499.75  # H90
$ws.Cells.Item(90, 10).Value = 76499.75  # J90
$ws.Cells.Item(90, 12).Value = 229499.25  # L90
$ws.Cells.Item(90, 14).Value = -241979.25  # N90
$ws.Cells.Item(93, 8).Value = 38999.4  # H93
$ws.Cells.Item(93, 10).Value = 38999.4  # J93
$ws.Cells.Item(93, 12).Value = 38999.4  # L93
$ws.Cells.Item(93, 14).Value = -42743.4  # N93
$ws.Cells.Item(102, 8).Value = 2058.4  # H102
$ws.Cells.Item(102, 9).Value = 1798.2812  # I102
$ws.Cells.Item(102, 10).Value = 4833  # J102
$ws.Cells.Item(102, 11).Value = 1798.2812  # K102
$ws.Cells.Item(102, 12).Value = 4833  # L102
$ws.Cells.Item(102, 13).Value = -176.2811999999999  # M102
$ws.Cells.Item(102, 14).Value = -8077  # N102
$ws.Cells.Item(122, 8).Value = 2967  # H122
$ws.Cells.Item(122, 9).Value = 2081.75  # I122
$ws.Cells.Item(122, 10).Value = 4737.5  # J122
$ws.Cells.Item(122, 11).Value = 6245.25  # K122
$ws.Cells.Item(122, 12).Value = 14212.5  # L122
$ws.Cells.Item(122, 13).Value = -3795.25  # M122
$ws.Cells.Item(122, 14).Value = -19112.5  # N122
$ws.Cells.Item(132, 8).Value = 2478.037  # H132
$ws.Cells.Item(132, 9).Value = 2454.4583  # I132
$ws.Cells.Item(132, 10).Value = 2666.6667  # J132
$ws.Cells.Item(132, 11).Value = 7363.374899999999  # K132
$ws.Cells.Item(132, 12).Value = 8000.000100000001  # L132
$ws.Cells.Item(132, 13).Value = -4833.374899999999  # M132
$ws.Cells.Item(132, 14).Value = -13060.0001  # N132

# Sheet 7 = LTW
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(7, 8).Value = 2948.8  # H7
$ws.Cells.Item(13, 8).Value = 6665  # H13
$ws.Cells.Item(13, 10).Value = 6665  # J13
$ws.Cells.Item(13, 12).Value = 6665  # L13
$ws.Cells.Item(13, 14).Value = -6945  # N13
$ws.Cells.Item(22, 8).Value = 1699.4286  # H22
$ws.Cells.Item(22, 9).Value = 798.5  # I22
$ws.Cells.Item(22, 11).Value = 798.5  # K22
$ws.Cells.Item(22, 13).Value = -503.5  # M22
$ws.Cells.Item(27, 8).Value = 1699.4286  # H27
$ws.Cells.Item(27, 9).Value = 798.5  # I27
$ws.Cells.Item(27, 11).Value = 798.5  # K27
$ws.Cells.Item(27, 13).Value = -691.5  # M27
$ws.Cells.Item(46, 8).Value = 2474.625  # H46
$ws.Cells.Item(46, 9).Value = 1099.6666  # I46
$ws.Cells.Item(46, 11).Value = 1099.6666  # K46
$ws.Cells.Item(46, 13).Value = -911.6666  # M46
$ws.Cells.Item(68, 8).Value = 2513.2  # H68
$ws.Cells.Item(68, 10).Value = 0  # J68
$ws.Cells.Item(68, 12).Value = 0  # L68
$ws.Cells.Item(68, 14).ClearContents()  # N68
$ws.Cells.Item(71, 8).Value = 2513.2  # H71
$ws.Cells.Item(71, 10).Value = 0  # J71
$ws.Cells.Item(71, 12).Value = 0  # L71
$ws.Cells.Item(71, 14).ClearContents()  # N71
$ws.Cells.Item(82, 8).Value = 1903.5  # H82
$ws.Cells.Item(82, 9).Value = 1722.55  # I82
$ws.Cells.Item(82, 10).Value = 2506.6667  # J82
$ws.Cells.Item(82, 11).Value = 1722.55  # K82
$ws.Cells.Item(82, 12).Value = 2506.6667  # L82
$ws.Cells.Item(82, 13).Value = -1361.55  # M82
$ws.Cells.Item(82, 14).Value = -3228.6667  # N82
$ws.Cells.Item(85, 8).Value = 1903.5  # H85
$ws.Cells.Item(85, 9).Value = 1722.55  # I85
$ws.Cells.Item(85, 10).Value = 2506.6667  # J85
$ws.Cells.Item(85, 11).Value = 1722.55  # K85
$ws.Cells.Item(85, 12).Value = 2506.6667  # L85
$ws.Cells.Item(85, 13).Value = -474.55  # M85
$ws.Cells.Item(85, 14).Value = -5002.6667  # N85
$ws.Cells.Item(126, 8).Value = 2948.8  # H126
$ws.Cells.Item(130, 8).Value = 0  # H130
$ws.Cells.Item(130, 10).Value = 0  # J130
$ws.Cells.Item(130, 12).Value = 0  # L130
$ws.Cells.Item(130, 14).ClearContents()  # N130
$ws.Cells.Item(132, 8).Value = 2107.4856  # H132
$ws.Cells.Item(132, 9).Value = 1264.4286  # I132
$ws.Cells.Item(132, 10).Value = 3372.0715  # J132
$ws.Cells.Item(132, 11).Value = 3793.2858  # K132
$ws.Cells.Item(132, 12).Value = 10116.2145  # L132
$ws.Cells.Item(132, 13).Value = -1263.2858  # M132
$ws.Cells.Item(132, 14).Value = -15176.2145  # N132
$ws.Cells.Item(134, 8).Value = 64500  # H134
$ws.Cells.Item(134, 10).Value = 64500  # J134
$ws.Cells.Item(134, 12).Value = 64500  # L134
$ws.Cells.Item(134, 14).Value = -74640  # N134
$ws.Cells.Item(141, 8).Value = 99881.6  # H141
$ws.Cells.Item(141, 10).Value = 99881.6  # J141
$ws.Cells.Item(141, 12).Value = 99881.6  # L141
$ws.Cells.Item(141, 14).Value = -110241.6  # N141

# Sheet 8 = WVR
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(29, 8).Value = 10374.75  # H29
$ws.Cells.Item(29, 9).Value = 6166.5  # I29
$ws.Cells.Item(29, 10).Value = 22999.5  # J29
$ws.Cells.Item(29, 11).Value = 6166.5  # K29
$ws.Cells.Item(29, 12).Value = 22999.5  # L29
$ws.Cells.Item(29, 13).Value = -5876.5  # M29
$ws.Cells.Item(29, 14).Value = -23579.5  # N29
$ws.Cells.Item(33, 8).Value = 6000  # H33
$ws.Cells.Item(33, 10).Value = 6000  # J33
$ws.Cells.Item(33, 12).Value = 6000  # L33
$ws.Cells.Item(33, 14).Value = -6500  # N33
$ws.Cells.Item(36, 8).Value = 6000  # H36
$ws.Cells.Item(36, 10).Value = 6000  # J36
$ws.Cells.Item(36, 12).Value = 6000  # L36
$ws.Cells.Item(36, 14).Value = -6500  # N36
$ws.Cells.Item(69, 8).Value = 45311.25  # H69
$ws.Cells.Item(69, 9).Value = 1246  # I69
$ws.Cells.Item(69, 10).Value = 59999.668  # J69
$ws.Cells.Item(69, 11).Value = 1246  # K69
$ws.Cells.Item(69, 12).Value = 59999.668  # L69
$ws.Cells.Item(69, 13).Value = -497  # M69
$ws.Cells.Item(69, 14).Value = -61497.668  # N69
$ws.Cells.Item(72, 8).Value = 45311.25  # H72
$ws.Cells.Item(72, 9).Value = 1246  # I72
$ws.Cells.Item(72, 10).Value = 59999.668  # J72
$ws.Cells.Item(72, 11).Value = 3738  # K72
$ws.Cells.Item(72, 12).Value = 179999.004  # L72
$ws.Cells.Item(72, 13).Value = 6  # M72
$ws.Cells.Item(72, 14).Value = -187487.004  # N72
$ws.Cells.Item(81, 8).Value = 11469.8  # H81
$ws.Cells.Item(81, 9).Value = 13624.75  # I81
$ws.Cells.Item(81, 11).Value = 27249.5  # K81
$ws.Cells.Item(81, 13).Value = -26188.5  # M81
$ws.Cells.Item(84, 8).Value = 11469.8  # H84
$ws.Cells.Item(84, 9).Value = 13624.75  # I84
$ws.Cells.Item(84, 11).Value = 136247.5  # K84
$ws.Cells.Item(84, 13).Value = -130943.5  # M84
$ws.Cells.Item(96, 8).Value = 2696.8572  # H96
$ws.Cells.Item(96, 10).Value = 3178.6  # J96
$ws.Cells.Item(96, 12).Value = 3178.6  # L96
$ws.Cells.Item(96, 14).Value = -5924.6  # N96
$ws.Cells.Item(132, 8).Value = 18521438  # H132
$ws.Cells.Item(132, 9).Value = 20002464  # I132
$ws.Cells.Item(132, 10).Value = 8624.75  # J132
$ws.Cells.Item(132, 11).Value = 60007392  # K132
$ws.Cells.Item(132, 12).Value = 25874.25  # L132
$ws.Cells.Item(132, 13).Value = -60004862  # M132
$ws.Cells.Item(132, 14).Value = -30934.25  # N132
